# Update the lattice-multiplication exercise table: 15 cells (5 rows x 3
# columns) get new "A x B" problems (and the matching digit breakdown
# lines) substituted in place of the old ones, per the target revision.
#
# Each table cell holds a single run with five text segments separated by
# <w:br/>:
#   1) "A x B"          - the problem statement
#   2) "  d1    d2"     - the two digits of B, spaced out (xml:space=preserve)
#   3) "  ----"         - separator (unchanged, xml:space=preserve)
#   4) "a1|    |"       - first digit of A, lattice row marker
#   5) "a2|    |"       - second digit of A, lattice row marker
#
# We rebuild each changed cell's Range via InsertXML with a full
# WordprocessingML package fragment so the serialized <w:t> runs retain
# their original xml:space="preserve" attribute (a plain Range.Text / Find
# & Replace assignment drops that attribute on rewritten text nodes).

function Set-LatticeCell {
    param(
        $Table,
        [int]$Row,
        [int]$Col,
        [string]$Problem,
        [string]$Factors,
        [string]$D1,
        [string]$D2
    )

    $cell = $Table.Cell($Row, $Col)
    $range = $cell.Range

    $xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
           '<?mso-application progid="Word.Document"?>' +
           '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
           '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
           '<pkg:xmlData>' +
           '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
           '<w:body>' +
           '<w:p><w:r><w:rPr><w:sz w:val="32"/></w:rPr>' +
           '<w:t>' + $Problem + '</w:t><w:br/>' +
           '<w:t xml:space="preserve">' + $Factors + '</w:t><w:br/>' +
           '<w:t xml:space="preserve">  ----</w:t><w:br/>' +
           '<w:t>' + $D1 + '</w:t><w:br/>' +
           '<w:t>' + $D2 + '</w:t>' +
           '</w:r></w:p>' +
           '</w:body></w:document>' +
           '</pkg:xmlData></pkg:part></pkg:package>'

    $range.InsertXML($xml) | Out-Null
}

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$cells = @(
    @{Row=1; Col=1; Problem="53 x 84"; Factors="  8    4"; D1="5|    |"; D2="3|    |"},
    @{Row=1; Col=2; Problem="38 x 16"; Factors="  1    6"; D1="3|    |"; D2="8|    |"},
    @{Row=1; Col=3; Problem="69 x 16"; Factors="  1    6"; D1="6|    |"; D2="9|    |"},
    @{Row=2; Col=1; Problem="25 x 93"; Factors="  9    3"; D1="2|    |"; D2="5|    |"},
    @{Row=2; Col=2; Problem="90 x 51"; Factors="  5    1"; D1="9|    |"; D2="0|    |"},
    @{Row=2; Col=3; Problem="99 x 89"; Factors="  8    9"; D1="9|    |"; D2="9|    |"},
    @{Row=3; Col=1; Problem="78 x 20"; Factors="  2    0"; D1="7|    |"; D2="8|    |"},
    @{Row=3; Col=2; Problem="97 x 56"; Factors="  5    6"; D1="9|    |"; D2="7|    |"},
    @{Row=3; Col=3; Problem="35 x 83"; Factors="  8    3"; D1="3|    |"; D2="5|    |"},
    @{Row=4; Col=1; Problem="33 x 43"; Factors="  4    3"; D1="3|    |"; D2="3|    |"},
    @{Row=4; Col=2; Problem="32 x 47"; Factors="  4    7"; D1="3|    |"; D2="2|    |"},
    @{Row=4; Col=3; Problem="86 x 10"; Factors="  1    0"; D1="8|    |"; D2="6|    |"},
    @{Row=5; Col=1; Problem="36 x 44"; Factors="  4    4"; D1="3|    |"; D2="6|    |"},
    @{Row=5; Col=2; Problem="54 x 99"; Factors="  9    9"; D1="5|    |"; D2="4|    |"},
    @{Row=5; Col=3; Problem="18 x 42"; Factors="  4    2"; D1="1|    |"; D2="8|    |"}
)

foreach ($item in $cells) {
    Set-LatticeCell $t $item.Row $item.Col $item.Problem $item.Factors $item.D1 $item.D2
}

Write-Host "Lattice multiplication table updated."
